# Auto-generated Excel COM-interop edit script
# Applies cell-level updates to Sheet1 per the crypto price/volume refresh diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.711.23"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.339.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -7.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.327.89"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.89%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.615"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.31"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.867.77"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.118"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.335.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.74"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.69"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "63.604.96"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.970"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.30"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.04"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.29"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.59"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.64"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.38%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.41"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "577.72"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.64"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.147"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.54"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -10.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.151.71"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.19"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0403"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.62"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.76"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.36%  "

Write-Host "Applied 102 cell updates."
